$d = $word.ActiveDocument

# =====================================================================
# 1) Locate the empty, right-aligned paragraph that immediately follows
#    the "...Rediger vare vinduet." caption and precedes the "Test"
#    heading. This is the paragraph that gets turned into the new
#    "Controllers" / "Cache" subsections.
# =====================================================================
$anchor = $d.Content
$anchor.Find.Execute("Rediger vare vinduet", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $anchor.Find.Found) { throw "anchor paragraph (Rediger vare vinduet) not found" }
$captionPar = $anchor.Paragraphs(1)
$cur = $captionPar.Next()

# --- 1a) Replace the empty paragraph with the "Controllers" heading (Overskrift4) ---
$curRange = $cur.Range
$curRange.Collapse(1)
$xml = @'
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:pStyle w:val="Overskrift4"/></w:pPr><w:proofErr w:type="spellStart"/><w:r><w:t>Controllers</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@
$curRange.InsertXML($xml)
$cur = $captionPar.Next()

# --- 1b) Insert the "Samlet set..." body paragraph right after it ---
$cur.Range.InsertParagraphAfter() | Out-Null
$cur = $cur.Next()
$curRange = $cur.Range
$curRange.Collapse(1)
$xml = @'
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t>Samlet set</w:t></w:r><w:r><w:t xml:space="preserve"> udgør </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>controllerne</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> det</w:t></w:r><w:r><w:t>,</w:t></w:r><w:r><w:t xml:space="preserve"> der svarer til</w:t></w:r><w:r><w:t xml:space="preserve"> busin</w:t></w:r><w:r><w:t>ess </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>logic</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>layer</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> fra </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Fridge-app</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>, og de varetager altså kommunikationen med den eksterne database</w:t></w:r><w:r><w:t xml:space="preserve">. </w:t></w:r><w:r><w:t xml:space="preserve">Hver </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>controller</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> har sit ansvar</w:t></w:r><w:r><w:t>.</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>LisView</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t>står for at indlæse og præsentere alle varer på en liste</w:t></w:r><w:r><w:t xml:space="preserve">, og giver mulighed for at slette en udvalgt vare. </w:t></w:r><w:r><w:t>Den står også for at omdirigere</w:t></w:r><w:r><w:t xml:space="preserve"> til både </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Additem</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>viewed</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> og </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>EditItem</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>viewed</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">. Inden den omdirigere til </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>EditItem</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">, finder den først den </w:t></w:r><w:r><w:t xml:space="preserve">varer, der skal sendes med til denne før der omdirigeres. </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>AddItem</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>controlleren</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> har til ansvar at tilføje varer til databasen</w:t></w:r><w:r><w:t xml:space="preserve">, hvis den finder en tilsvarende vare på den nuværende liste sørge den for blot at tælle antallet af eksisterende vare op i stedet for at indsætte en </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>duplet</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">. </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>EditItem</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>controlleren</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> står for redigering</w:t></w:r><w:r><w:t xml:space="preserve"> og opdatering af en valgt vare. </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>EditItem</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> sørger også for at </w:t></w:r><w:r><w:t xml:space="preserve">en vare ikke kan opdateres så den bliver en </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>duplet</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> på listen, dette håndteres ligesom i</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>AddItem</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">. </w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@
$curRange.InsertXML($xml)

# --- 1c) Insert the "Cache" heading (Overskrift4) right after it ---
$cur.Range.InsertParagraphAfter() | Out-Null
$cur = $cur.Next()
$curRange = $cur.Range
$curRange.Collapse(1)
$xml = @'
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:pStyle w:val="Overskrift4"/></w:pPr><w:proofErr w:type="spellStart"/><w:r><w:t>Cache</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@
$curRange.InsertXML($xml)

# --- 1d) Insert the Cache body paragraph right after it ---
$cur.Range.InsertParagraphAfter() | Out-Null
$cur = $cur.Next()
$curRange = $cur.Range
$curRange.Collapse(1)
$xml = @'
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t xml:space="preserve">Til web applikationen er der blevet implementeret en </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>cache</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">. Denne </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>cache</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> indeholder alle de ting, som er fælles for</w:t></w:r><w:r><w:t xml:space="preserve"> alle </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>controllers</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">. Klassen består kun </w:t></w:r><w:r><w:t>af</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>properties</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>, der</w:t></w:r><w:r><w:t xml:space="preserve"> alle er </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>stati</w:t></w:r><w:r><w:t>c</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">. Begrundelsen for dette er netop at det er fælles data, som alle </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>controllers</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> skal arbejde på. Det er også i </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>cachen</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> at facaden til data </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>acce</w:t></w:r><w:r><w:t>ss</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>layer</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> ligger, dette er fordi der på alle tidspunkter kun må være en facade. </w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@
$curRange.InsertXML($xml)

# =====================================================================
# 2) In the "Test" section, split the run that talks about testing for
#    an existing item so a lastRenderedPageBreak sits right before
#    "eksistere et item, ...".
# =====================================================================
$rng2 = $d.Content
$rng2.Find.Execute("om der eksistere", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $rng2.Find.Found) { throw "anchor text (om der eksistere) not found" }
$splitPar = $rng2.Paragraphs(1)
$splitParRange = $splitPar.Range
$splitParRange.Collapse(1)
$xml = @'
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p w:rsidR="00E7168B" w:rsidRDefault="001C623A" w:rsidP="00E7168B"><w:r><w:t xml:space="preserve">Et af problemerne ved den måde, som funktionerne i applikationen er opbygget på, er at mange af dem returnere et </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>View</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> eller et </w:t></w:r><w:r w:rsidR="008D387C"><w:t>’</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>RedirectToAction</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidR="008D387C"><w:t>’</w:t></w:r><w:r w:rsidR="00910F0D"><w:t xml:space="preserve"> (</w:t></w:r><w:r w:rsidR="00976C85"><w:fldChar w:fldCharType="begin"/></w:r><w:r w:rsidR="00910F0D"><w:instrText xml:space="preserve"> REF _Ref419899941 \h </w:instrText></w:r><w:r w:rsidR="00976C85"><w:fldChar w:fldCharType="separate"/></w:r><w:r w:rsidR="00910F0D"><w:t xml:space="preserve">Figur </w:t></w:r><w:r w:rsidR="00910F0D"><w:rPr><w:noProof/></w:rPr><w:t>1</w:t></w:r><w:r w:rsidR="00976C85"><w:fldChar w:fldCharType="end"/></w:r><w:r w:rsidR="00910F0D"><w:t>)</w:t></w:r><w:r><w:t xml:space="preserve"> efter et item har gennemgået anden logik i samme funktion. Et eksempel på dette kunne være i </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>EditItemController</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> hvor funktionen </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>UpdateItem</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> er implementeret. Denne funktion tager imod en </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>FormCollection</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> hvori </w:t></w:r><w:r w:rsidR="008D387C"><w:t xml:space="preserve">de nye værdier der skal ændres i det gamle item ligger i. logikken i funktionen ændre det gamle item til de nye værdier, men returnere det </w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidR="008D387C"><w:t>view</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidR="008D387C"><w:t xml:space="preserve">, som skal vises efter </w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidR="008D387C"><w:t>itemet</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidR="008D387C"><w:t xml:space="preserve"> er ændret. Dvs. det ikke er muligt bare at sammenligne returværdien fra funktionen med det der forventes at det pågældende item er ændret til. I stedet skal der testes på om der </w:t></w:r><w:r w:rsidR="008D387C"><w:lastRenderedPageBreak/><w:t>eksistere et item, i listen med items, der har de værdier der skulle ændres. Dette går igen med alle de funktioner der har ’</w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidR="008D387C"><w:t>ActionResult</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidR="008D387C"><w:t>’ som retur parameter.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@
$splitParRange.InsertXML($xml)

Write-Host "Edit complete. Paragraph count:" $d.Paragraphs.Count
